$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the data columns (Price, Volume) keep their original text storage type,
# matching the source workbook where every D/E cell is stored as text (inline string).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "57.373.84"
$ws.Range("E2").Value = "  -4.90%  "
$ws.Range("D3").Value = "2.907.30"
$ws.Range("E3").Value = "  -2.16%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "548.35"
$ws.Range("E5").Value = "  -1.75%  "
$ws.Range("D6").Value = "125.29"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "2.908.46"
$ws.Range("E8").Value = "  -2.08%  "
$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D9").Value = "0.502"
$ws.Range("E9").Value = "  +1.95%  "
$ws.Range("D10").Value = "0.123"
$ws.Range("E10").Value = "  -7.21%  "
$ws.Range("D11").Value = "4.69"
$ws.Range("E11").Value = "  -8.65%  "
$ws.Range("D12").Value = "0.436"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "0.0000211"
$ws.Range("E13").Value = "  -4.58%  "
$ws.Range("D14").Value = "32.30"
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").Value = "3.382.39"
$ws.Range("E16").Value = "  -2.19%  "
$ws.Range("D17").Value = "2.895.44"
$ws.Range("E17").Value = "  -2.33%  "
$ws.Range("D18").Value = "6.53"
$ws.Range("E18").Value = "  +6.89%  "
$ws.Range("D19").Value = "57.260.39"
$ws.Range("E19").Value = "  -5.39%  "
$ws.Range("D20").Value = "406.32"
$ws.Range("E20").Value = "  -5.33%  "
$ws.Range("D21").Value = "12.82"
$ws.Range("E21").Value = "  -1.13%  "
$ws.Range("D22").Value = "0.672"
$ws.Range("E22").Value = "  +2.61%  "
$ws.Range("D23").Value = "6.83"
$ws.Range("E23").Value = "  -3.64%  "
$ws.Range("D24").Value = "12.67"
$ws.Range("E24").Value = "  -1.43%  "
$ws.Range("D25").Value = "77.60"
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").Value = "0.997"
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("D28").Value = "2.45"
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").Value = "7.22"
$ws.Range("E29").Value = "  +1.92%  "
$ws.Range("D30").Value = "1.93"
$ws.Range("E30").Value = "  +2.89%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "5.97"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "24.67"
$ws.Range("E32").Value = "  -1.77%  "
$ws.Range("D33").Value = "0.0988"
$ws.Range("E33").Value = "  +6.41%  "
$ws.Range("D34").Value = "0.918"
$ws.Range("E34").Value = "  -2.98%  "
$ws.Range("D35").Value = "5.43"
$ws.Range("E35").Value = "  -1.02%  "
$ws.Range("D36").Value = "2.01"
$ws.Range("E36").Value = "  -10.65%  "
$ws.Range("D37").Value = "48.20"
$ws.Range("E37").Value = "  -2.06%  "
$ws.Range("D38").Value = "8.26"
$ws.Range("E38").Value = "  +6.98%  "
$ws.Range("D39").Value = "0.0₃0631"
$ws.Range("E39").Value = "  -3.78%  "
$ws.Range("D40").Value = "0.106"
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("D41").Value = "0.0338"
$ws.Range("E41").Value = "  -4.68%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "2.43"
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.620.25"
$ws.Range("E43").Value = "  -1.50%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").Value = "361.19"
$ws.Range("E44").Value = "  -3.02%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "119.82"
$ws.Range("E47").Value = "  -2.06%  "
$ws.Range("E48").Value = "  +1.82%  "
$ws.Range("D49").Value = "1.93"
$ws.Range("E49").Value = "  -0.60%  "
$ws.Range("D50").Value = "22.85"
$ws.Range("E50").Value = "  -1.34%  "
$ws.Range("D51").Value = "1.95"
$ws.Range("E51").Value = "  -2.40%  "
